$d = $word.ActiveDocument

$pairs = @(
    @("48×21=", "24×83="),
    @("74×96=", "40×71="),
    @("58×87=", "82×73="),
    @("80×87=", "79×84="),
    @("13×65=", "31×47="),
    @("33×33=", "13×17="),
    @("11×62=", "15×18="),
    @("35×20=", "13×94="),
    @("62×70=", "91×70="),
    @("46×78=", "65×17="),
    @("84×85=", "72×61="),
    @("69×57=", "67×61="),
    @("90×49=", "29×97="),
    @("91×69=", "63×60="),
    @("69×73=", "17×80="),
    @("79×22=", "51×73="),
    @("56×32=", "91×84="),
    @("26×62=", "59×42="),
    @("52×21=", "35×17="),
    @("49×84=", "88×35="),
    @("51×37=", "45×53="),
    @("75×81=", "91×30="),
    @("97×22=", "70×58="),
    @("20×41=", "39×81="),
    @("88×60=", "88×57=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
